$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B10 value: 175 -> 174
$ws.Range("B10").Value = 174

# Add new label in D10 with bold red font style
$ws.Range("D10").Value = "CASO VIGA BIEN DIMENSIONADA"
$ws.Range("D10").Font.Bold = $true
$ws.Range("D10").Font.Color = 3355647
$ws.Range("D10").Font.Name = "Arial"
$ws.Range("D10").Font.Size = 10

# Update selection to D10
$ws.Range("D10").Select() | Out-Null
